# Apply symbol-list / price updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.98'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.67'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.399'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05762'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.435'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.335'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8117'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9080'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1448'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07352'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03145'
$ws.Range("E12").Value = '11BitrueCoinBTRBestin24h'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02973'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09419'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.935'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001595'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04822'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005841'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006154'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'UpBots'
$ws.Range("C20").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.007491'
$ws.Range("E20").Value = '19UpBotsUBXTWorstin24h'
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004066'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("B22").Value = 'BitKan'
$ws.Range("C22").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0009948'
$ws.Range("E22").Value = '21BitKanKAN'
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0001501'
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.749'
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.199'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.3280'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("B27").Value = 'ProBitToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1330'
$ws.Range("E27").Value = '26ProBitTokenPROB'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03903'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006794'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1074'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002411'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007329'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005639'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1676'
